$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '66.244.19'
$c.Style = $origStyle
$ws.Range('E2').Value = '  -1.25%  '

$c = $ws.Range('D3')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '3.514.93'
$c.Style = $origStyle
$ws.Range('E3').Value = '  -0.35%  '

$ws.Range('E4').Value = '  +0.04%  '

$c = $ws.Range('D5')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '585.96'
$c.Style = $origStyle
$ws.Range('E5').Value = '  +6.30%  '

$c = $ws.Range('D6')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '179.21'
$c.Style = $origStyle
$ws.Range('E6').Value = '  -5.28%  '

$c = $ws.Range('D7')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.635'
$c.Style = $origStyle
$ws.Range('E7').Value = '  +4.78%  '

$ws.Range('E8').Value = '  +0.04%  '

$ws.Range('E9').Value = '  +1.37%  '

$c = $ws.Range('D10')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.164'
$c.Style = $origStyle
$ws.Range('E10').Value = '  +5.56%  '

$c = $ws.Range('D11')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '55.79'
$c.Style = $origStyle
$ws.Range('E11').Value = '  +2.09%  '

$ws.Range('E12').Value = '  +3.70%  '

$ws.Range('E13').Value = '  -0.60%  '

$c = $ws.Range('D14')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '4.074.36'
$c.Style = $origStyle
$ws.Range('E14').Value = '  -0.08%  '

$c = $ws.Range('D15')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '3.516.62'
$c.Style = $origStyle
$ws.Range('E15').Value = '  -0.07%  '

$ws.Range('E16').Value = '  +0.21%  '

$c = $ws.Range('D17')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '18.42'
$c.Style = $origStyle
$ws.Range('E17').Value = '  +1.31%  '

$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$c = $ws.Range('D18')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '66.226.02'
$c.Style = $origStyle
$ws.Range('E18').Value = '  -1.29%  '

$ws.Range('B19').Value = 'Uniswap'
$ws.Range('C19').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$c = $ws.Range('D19')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '12.09'
$c.Style = $origStyle
$ws.Range('E19').Value = '  +1.43%  '

$ws.Range('E20').Value = '  +2.29%  '

$c = $ws.Range('D21')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '415.83'
$c.Style = $origStyle
$ws.Range('E21').Value = '  -3.99%  '

$ws.Range('E22').Value = '  +10.91%  '

$c = $ws.Range('D23')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '4.45'
$c.Style = $origStyle
$ws.Range('E23').Value = '  +7.22%  '

$c = $ws.Range('D24')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '85.24'
$c.Style = $origStyle
$ws.Range('E24').Value = '  +0.02%  '

$c = $ws.Range('D25')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '13.58'
$c.Style = $origStyle
$ws.Range('E25').Value = '  +13.08%  '

$c = $ws.Range('D26')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '11.13'
$c.Style = $origStyle
$ws.Range('E26').Value = '  +0.22%  '

$c = $ws.Range('D27')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '2.87'
$c.Style = $origStyle
$ws.Range('E27').Value = '  -1.11%  '

$c = $ws.Range('D28')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '6.05'
$c.Style = $origStyle
$ws.Range('E28').Value = '  -1.62%  '

$c = $ws.Range('D29')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '9.20'
$c.Style = $origStyle
$ws.Range('E29').Value = '  +2.43%  '

$c = $ws.Range('D30')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '30.45'
$c.Style = $origStyle
$ws.Range('E30').Value = '  +0.66%  '

$c = $ws.Range('D31')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '6.67'
$c.Style = $origStyle
$ws.Range('E31').Value = '  +0.01%  '

$c = $ws.Range('D32')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '607.50'
$c.Style = $origStyle
$ws.Range('E32').Value = '  -5.53%  '

$ws.Range('E33').Value = '  +0.90%  '

$ws.Range('E34').Value = '  +0.63%  '

$c = $ws.Range('D35')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '60.53'
$c.Style = $origStyle
$ws.Range('E35').Value = '  +1.47%  '

$c = $ws.Range('D36')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.154'
$c.Style = $origStyle
$ws.Range('E36').Value = '  +7.46%  '

$ws.Range('B37').Value = 'PEPE'
$ws.Range('C37').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$c = $ws.Range('D37')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.0₃0804'
$c.Style = $origStyle
$ws.Range('E37').Value = '  -2.04%  '

$ws.Range('B38').Value = 'Dai'
$ws.Range('C38').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$c = $ws.Range('D38')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = $origStyle
$ws.Range('E38').Value = '  +0.10%  '

$ws.Range('B39').Value = 'Stacks'
$ws.Range('C39').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$c = $ws.Range('D39')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '3.67'
$c.Style = $origStyle
$ws.Range('E39').Value = '  +9.37%  '

$c = $ws.Range('D40')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '37.00'
$c.Style = $origStyle
$ws.Range('E40').Value = '  -3.95%  '

$c = $ws.Range('D41')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.386'
$c.Style = $origStyle
$ws.Range('E41').Value = '  -1.02%  '

$c = $ws.Range('D42')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '3.260.09'
$c.Style = $origStyle
$ws.Range('E42').Value = '  +7.75%  '

$c = $ws.Range('D43')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.Style = $origStyle
$ws.Range('E43').Value = '  +0.06%  '

$c = $ws.Range('D44')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '2.99'
$c.Style = $origStyle
$ws.Range('E44').Value = '  +4.45%  '

$c = $ws.Range('D45')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '3.34'
$c.Style = $origStyle
$ws.Range('E45').Value = '  +0.14%  '

$ws.Range('E46').Value = '  -2.86%  '

$ws.Range('E47').Value = '  +1.32%  '

$ws.Range('B48').Value = 'Stellar'
$ws.Range('C48').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$c = $ws.Range('D48')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.134'
$c.Style = $origStyle
$ws.Range('E48').Value = '  +2.43%  '

$ws.Range('B49').Value = 'WEMIXToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$c = $ws.Range('D49')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '2.69'
$c.Style = $origStyle
$ws.Range('E49').Value = '  -6.01%  '

$c = $ws.Range('D50')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '8.67'
$c.Style = $origStyle
$ws.Range('E50').Value = '  -0.05%  '

$c = $ws.Range('D51')
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '140.04'
$c.Style = $origStyle
$ws.Range('E51').Value = '  -1.63%  '
